$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.484.36"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "2.191.99"
$ws.Range("E3").Value = "  -7.28%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'483.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").Value = "'124.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.83%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "'0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.59%  "
$ws.Range("D9").Value = "2.210.45"
$ws.Range("E9").Value = "  -6.64%  "
$ws.Range("D10").Value = "'0.0912"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.94%  "
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'4.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.90%  "
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.312"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "2.580.32"
$ws.Range("E14").Value = "  -7.27%  "
$ws.Range("D15").Value = "'20.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "53.401.08"
$ws.Range("E16").Value = "  -4.67%  "
$ws.Range("E17").Value = "  -3.78%  "
$ws.Range("D18").Value = "2.189.89"
$ws.Range("E18").Value = "  -6.88%  "
$ws.Range("D19").Value = "'9.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.69%  "
$ws.Range("D20").Value = "'3.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "'294.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'62.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").Value = "'0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").Value = "'0.364"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").Value = "2.289.02"
$ws.Range("E27").Value = "  -7.43%  "
$ws.Range("D28").Value = "'0.144"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "'6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").Value = "'165.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.32%  "
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'0.994"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "0.0₃0657"
$ws.Range("E34").Value = "  -7.27%  "
$ws.Range("D35").Value = "'5.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'17.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("D39").Value = "'0.819"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.89%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'35.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").Value = "'0.365"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D43").Value = "'1.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.81%  "
$ws.Range("D44").Value = "'3.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.68%  "
$ws.Range("D45").Value = "'123.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("D46").Value = "'4.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").Value = "'0.530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.74%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'228.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D50").Value = "'0.0466"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("E51").Value = "  -4.02%  "
